$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the skill IDs in column A to their new upper-case spellings
# (write the SKILL* rows first, then the NORMALATTACK* rows, so the shared
# string table is appended in that order) ---
$ws.Range("A6").Value = "SKILL1"
$ws.Range("A7").Value = "SKILL2"
$ws.Range("A8").Value = "SKILL3"
$ws.Range("A9").Value = "SKILL4"

$ws.Range("A2").Value = "NORMALATTACK1"
$ws.Range("A3").Value = "NORMALATTACK2"
$ws.Range("A4").Value = "NORMALATTACK3"
$ws.Range("A5").Value = "NORMALTHUMP"

# --- Every row's NextLevelID (column C) now points at NORMALATTACK2 ---
$ws.Range("C2:C9").Value = "NORMALATTACK2"

# --- Clear the leftover banded fill/border from C4:C9 so they match the
# plain text style already used by C2:C3 (copy C2's format onto them) ---
$ws.Range("C2").Copy()
$ws.Range("C4:C9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Bump the AnimaState (column G) ids for the new Skill rows ---
$ws.Range("G6").Value = 101
$ws.Range("G7").Value = 102
$ws.Range("G8").Value = 103
$ws.Range("G9").Value = 104

# --- Apply a text number format to the NextLevelID table column (matches
# the dxf/dataDxfId the author recorded for that column in the table) ---
$ws.Range("C2:C9").NumberFormat = "@"

# --- Restore the selection the author left behind: C2:C5 then C6:C9 added
# to it, with C6 the active cell ---
$r1 = $ws.Range("C2:C5")
$r2 = $ws.Range("C6:C9")
$excel.Union($r1, $r2).Select() | Out-Null
